$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new JSON entry shifted the time-bucket ordering, which swapped the
# data that lives in row 2 vs row 3 (title, timestamp, historical
# distance, uri) while the "time bucket" column (D) stays the same.

$titleA = $ws.Range("A2").Value2
$titleB = $ws.Range("A3").Value2
$ws.Range("A2").Value = $titleB
$ws.Range("A3").Value = $titleA

$timeA = $ws.Range("B2").Value2
$timeB = $ws.Range("B3").Value2
$ws.Range("B2").Value = $timeB
$ws.Range("B3").Value = $timeA

$distA = $ws.Range("C2").Value2
$distB = $ws.Range("C3").Value2
$ws.Range("C2").Value = $distB
$ws.Range("C3").Value = $distA

$uriA = $ws.Range("E2").Value2
$uriB = $ws.Range("E3").Value2
$ws.Range("E2").Value = $uriB
$ws.Range("E3").Value = $uriA
